# Pub_Style.xlsx edit: add "Year" column (D) populated for dissertation
# (unpublished) references, and refocus the AutoFilter from column A
# (Full Reference code) onto column C ("Dissertation (unpub)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the existing AutoFilter criteria (column A: CC06AB / CC06AC)
#    so that all rows are first shown, before we touch hidden rows --
#    writing into a currently-hidden row causes Excel to stamp a bogus
#    explicit row height, so we reset filtering before editing cells.
$ws.AutoFilterMode = $false

# 2) Re-apply the AutoFilter on the same A1:C478 range, this time
#    filtering column C (index 3, 1-based) for "Dissertation (unpub)".
#    This both rewrites the <autoFilter> definition and recalculates
#    every row's hidden state to match the new criteria.
$ws.Range("A1:C478").AutoFilter(3, @("Dissertation (unpub)"), 7)

# 3) Add the new "Year" header and per-row publication years for every
#    dissertation (unpublished) row. These rows are now visible (not
#    hidden) thanks to step 2, so no stray row-height gets introduced.
$ws.Range("D1").Value = "Year"

$ws.Range("D2").Value = 1992
$ws.Range("D65").Value = 1987
$ws.Range("D66").Value = 1987
$ws.Range("D67").Value = 1987
$ws.Range("D68").Value = 2007
$ws.Range("D100").Value = 2006
$ws.Range("D101").Value = 2006
$ws.Range("D102").Value = 2006
$ws.Range("D103").Value = 2006
$ws.Range("D114").Value = 2000
$ws.Range("D115").Value = 2000
$ws.Range("D116").Value = 2000
$ws.Range("D117").Value = 1987
$ws.Range("D118").Value = 1987
$ws.Range("D138").Value = 1995
$ws.Range("D139").Value = 1995
$ws.Range("D140").Value = 1995
$ws.Range("D141").Value = 1995
$ws.Range("D156").Value = 2005
$ws.Range("D168").Value = 1983
$ws.Range("D169").Value = 1983
$ws.Range("D178").Value = 2006
$ws.Range("D179").Value = 2006
$ws.Range("D221").Value = 1994
$ws.Range("D222").Value = 1994
$ws.Range("D250").Value = 1999
$ws.Range("D263").Value = 2002
$ws.Range("D283").Value = 2003
$ws.Range("D284").Value = 2003
$ws.Range("D285").Value = 2003
$ws.Range("D286").Value = 2003
$ws.Range("D298").Value = 2007
$ws.Range("D299").Value = 2007
$ws.Range("D300").Value = 2007
$ws.Range("D331").Value = 2008
$ws.Range("D333").Value = 2008
$ws.Range("D334").Value = 2006
$ws.Range("D335").Value = 2006
$ws.Range("D427").Value = 1984
$ws.Range("D428").Value = 1984
$ws.Range("D429").Value = 1984

# 4) Move the active selection to D428, matching where editing finished.
$ws.Range("D428").Select()
